$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.558139091228547
$ws.Range("C2").Value = 1.807629892617897
$ws.Range("D2").Value = 1.550484622924876
$ws.Range("E2").Value = 1.414055654904187
$ws.Range("B3").Value = 2.571058240846794
$ws.Range("C3").Value = 1.817304492547184
$ws.Range("D3").Value = 1.55689358682896
$ws.Range("E3").Value = 1.418476688626426
$ws.Range("B4").Value = 2.524033221911411
$ws.Range("C4").Value = 1.782010133492829
$ws.Range("D4").Value = 1.527468092681718
$ws.Range("E4").Value = 1.393836101180278
$ws.Range("B5").Value = 2.564825141152768
$ws.Range("C5").Value = 1.812988199412726
$ws.Range("D5").Value = 1.537805340834631
$ws.Range("E5").Value = 1.414522073045725
$ws.Range("B6").Value = 2.574585052914157
$ws.Range("C6").Value = 1.819960065524298
$ws.Range("D6").Value = 1.543419494542331
$ws.Range("E6").Value = 1.419532925750195
$ws.Range("B7").Value = 2.549393159961533
$ws.Range("C7").Value = 1.800608292030371
$ws.Range("D7").Value = 1.54268538216686
$ws.Range("E7").Value = 1.404978498878209
$ws.Range("B8").Value = 2.562226646976647
$ws.Range("C8").Value = 1.811121680113447
$ws.Range("D8").Value = 1.565222802774
$ws.Range("E8").Value = 1.409312192220709
$ws.Range("B9").Value = 2.569237339612119
$ws.Range("C9").Value = 1.816956055849428
$ws.Range("D9").Value = 1.557877905496985
$ws.Range("E9").Value = 1.419164020447879
$ws.Range("B10").Value = 2.259907216303427
$ws.Range("C10").Value = 1.584180934600336
$ws.Range("D10").Value = 1.375272804597484
$ws.Range("E10").Value = 1.27224744976302
$ws.Range("B11").Value = 2.194203204459443
$ws.Range("C11").Value = 1.532756793977353
$ws.Range("D11").Value = 1.329380741006684
$ws.Range("E11").Value = 1.232365762670961
$ws.Range("B12").Value = 1.865783552065658
$ws.Range("C12").Value = 1.280968368999795
$ws.Range("D12").Value = 1.082554188446728
$ws.Range("E12").Value = 0.9993215601058449
$ws.Range("B13").Value = 2.237924815214063
$ws.Range("C13").Value = 1.565073904619898
$ws.Range("D13").Value = 1.350229857885922
$ws.Range("E13").Value = 1.24896039519903
